# Rename the "Field_of_Education" sheet to "Education" and make it the
# active/selected tab (replacing "Enrolments" as the active sheet).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Field_of_Education")
$ws.Name = "Education"

$ws.Activate()
$ws.Select()
